$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("flussi_previsti")
$ws.Activate()

# Delete row 91 (regionali 2024 / astensione-only row) and shift the rest up.
$ws.Rows.Item(91).Delete()

# Reset the view: clear the "frozen" top-left scroll position and selection.
$ws.Range("C3").Select()
